# Apply "Propellor" test-case updates + refreshed correction-factor values
# to the ADDA/Linux/CorrectionFactors.xlsx workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Refresh existing correction-factor results (higher-precision reruns) ---
$ws.Range("E4").Value = 0.308310574
$ws.Range("E11").Value = 0.1170516473
$ws.Range("E12").Value = 0.3141289631
$ws.Range("E15").Value = 0.1324294717

# --- New "Propellor" section (rows 18-23) ---

# Section header row (bold label + note), mirrors the Sphere/Cone headers above
$ws.Range("A18").Value = "Propellor"
$ws.Range("A18").Font.Bold = $true
$ws.Range("B18").Value = "Considering F_x,F_y,F_z"

# Default 15 dpl test row
$ws.Range("A19").Value = 15
$ws.Range("B19").Value = 0.9
$ws.Range("C19").Value = 1.1859519223999999
$ws.Range("D19").Value = 2
$ws.Range("G19").Value = "Polystyrene propellor in water (Width 2, Height 2 micro m)"

# Default 15 dpl, dpl based on Lambda/RI
$ws.Range("A20").Value = 15
$ws.Range("B20").Value = 1
$ws.Range("C20").Value = 1.5
$ws.Range("D20").Value = 2
$ws.Range("G20").Value = "Default 15 dpl tests"

# Default 30 dpl test row
$ws.Range("A21").Value = 30
$ws.Range("B21").Value = 1
$ws.Range("C21").Value = 1.5
$ws.Range("D21").Value = 2
$ws.Range("G21").Value = "Default 30 dpl tests"

# Final results row (row 22 intentionally left blank, matching the gap
# used before every other "Final Results" row in this sheet)
$ws.Range("A23").Value = 15
$ws.Range("B23").Value = 1.0640000000000001
$ws.Range("C23").Value = 1.1859519223999999
$ws.Range("D23").Value = 2
$ws.Range("G23").Value = "Final Results for Polystyrene propellor in water (Width 2, Height 2 micro m)"

# Update dimension / current selection to reflect the new used range
$ws.Range("E6").Select()
